# Apply edits to worksheets '展览' (sheet index 1) and '全部类型' (sheet index 4)
$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)
foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Row 2
    $ws.Range("B2").Value = '2024-01-27'
    $ws.Range("C2").Formula = ""
    $ws.Range("D2").Value = '临泉路88号板桥里墨园E区1号省羽体中心 省羽体super速搏羽毛球馆'
    $ws.Range("E2").Value = '2024.01.27 10:00-01.28 17:00'
    $ws.Range("F2").Value = 1431
    $ws.Range("G2").Value = '不可售'
    $ws.Range("H2").Value = $false
    $ws.Range("I2").Value = 'https://show.bilibili.com/platform/detail.html?id=78076'
    $ws.Range("J2").Value = '//i1.hdslb.com/bfs/openplatform/202311/2v00jbxM1698999146733.jpeg'

    # Row 3
    $ws.Range("B3").Value = '2024-01-27'
    $ws.Range("C3").Value = '合肥·新春AG动漫游戏盛典热血plus'
    $ws.Range("D3").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E3").Value = '2024.01.27 10:00-01.28 17:30'
    $ws.Range("F3").Value = 7475
    $ws.Range("G3").Value = '65'
    $ws.Range("H3").Value = $true
    $ws.Range("I3").Value = 'https://show.bilibili.com/platform/detail.html?id=79566'
    $ws.Range("J3").Value = '//i1.hdslb.com/bfs/openplatform/202312/iJ1Dnmla1702029064983.jpeg'

    # Row 4
    $ws.Range("B4").Value = '2024-01-28'
    $ws.Range("C4").Value = '合肥·SP同人展·次元派对（取消）'
    $ws.Range("D4").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E4").Value = '2024.01.28 10:00-01.28 17:00'
    $ws.Range("F4").Value = 530
    $ws.Range("G4").Value = '已售罄'
    $ws.Range("H4").Value = $false
    $ws.Range("I4").Value = 'https://show.bilibili.com/platform/detail.html?id=80254'
    $ws.Range("J4").Value = '//i0.hdslb.com/bfs/openplatform/202312/9ClQwbVE1703668101900.jpeg'

    # Row 5
    $ws.Range("B5").Value = '2024-01-28'
    $ws.Range("C5").Value = '合肥·梦时空SPO1动漫展'
    $ws.Range("D5").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E5").Value = '2024.01.28 10:00-01.28 17:00'
    $ws.Range("F5").Value = 319
    $ws.Range("G5").Value = '258'
    $ws.Range("H5").Value = $false
    $ws.Range("I5").Value = 'https://show.bilibili.com/platform/detail.html?id=80246'
    $ws.Range("J5").Value = '//i0.hdslb.com/bfs/openplatform/202312/aHzqArm61703662347629.jpeg'

    # Row 6
    $ws.Range("B6").Value = '2024-01-28'
    $ws.Range("C6").Value = '肥东· 原神&崩铁&崩坏only'
    $ws.Range("D6").Value = '团结东路7号 巢湖宾馆'
    $ws.Range("E6").Value = '2024.01.28 10:00-01.28 17:00'
    $ws.Range("F6").Value = 22
    $ws.Range("G6").Value = '55'
    $ws.Range("H6").Value = $false
    $ws.Range("I6").Value = 'https://show.bilibili.com/platform/detail.html?id=80917'
    $ws.Range("J6").Value = '//i0.hdslb.com/bfs/openplatform/202401/UekMeUjQ1705462868391.jpeg'

    # Row 7
    $ws.Range("B7").Value = '2024-01-29'
    $ws.Range("C7").Value = '巢湖·原×铁×崩only'
    $ws.Range("D7").Value = '长江东路徽商城2幢B座(祥和地铁站C口步行370米) 曼斯顿尚品酒店'
    $ws.Range("E7").Value = '2024.01.29 10:00-01.29 17:00'
    $ws.Range("F7").Value = 14
    $ws.Range("G7").Value = '55'
    $ws.Range("H7").Value = $false
    $ws.Range("I7").Value = 'https://show.bilibili.com/platform/detail.html?id=80919'
    $ws.Range("J7").Value = '//i0.hdslb.com/bfs/openplatform/202401/9XumHIT31705464002179.jpeg'

    # Row 8
    $ws.Range("B8").Value = '2024-01-31'
    $ws.Range("C8").Value = '巢湖·原神&崩铁&崩坏only'
    $ws.Range("D8").Value = '仙满楼·麦肯希酒店 仙满楼·麦肯希酒店'
    $ws.Range("E8").Value = '2024.01.31 10:00-01.31 17:00'
    $ws.Range("F8").Value = 20
    $ws.Range("G8").Value = '55'
    $ws.Range("H8").Value = $false
    $ws.Range("I8").Value = 'https://show.bilibili.com/platform/detail.html?id=80944'
    $ws.Range("J8").Value = '//i0.hdslb.com/bfs/openplatform/202401/euD63Mlp1705479140627.jpeg'

    # Row 9
    $ws.Range("B9").Value = '2024-02-03'
    $ws.Range("C9").Value = '合肥·2024运动新春动漫庆典（全ip）'
    $ws.Range("D9").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E9").Value = '2024.02.03 09:30-02.04 17:00'
    $ws.Range("F9").Value = 5514
    $ws.Range("G9").Value = '65'
    $ws.Range("H9").Value = $true
    $ws.Range("I9").Value = 'https://show.bilibili.com/platform/detail.html?id=79963'
    $ws.Range("J9").Value = '//i0.hdslb.com/bfs/openplatform/202312/tBk3WVyX1702968658234.jpeg'

    # Row 10
    $ws.Range("B10").Value = '2024-02-04'
    $ws.Range("C10").Value = '合肥·六安lovelive only'
    $ws.Range("D10").Value = '健康东路7号 巢湖国际饭店'
    $ws.Range("E10").Value = '2024.02.04 10:00-02.04 17:00'
    $ws.Range("F10").Value = 9
    $ws.Range("G10").Value = '60'
    $ws.Range("H10").Value = $false
    $ws.Range("I10").Value = 'https://show.bilibili.com/platform/detail.html?id=80974'
    $ws.Range("J10").Value = '//i0.hdslb.com/bfs/openplatform/202401/wVVrdShB1705487994232.jpeg'

    # Row 11
    $ws.Range("B11").Value = '2024-02-04'
    $ws.Range("C11").Value = '合肥·第十二届次元之门动漫游戏博览会-赵乾景专场'
    $ws.Range("D11").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
    $ws.Range("E11").Value = '2024.02.04 11:30-02.04 17:00'
    $ws.Range("F11").Value = 142
    $ws.Range("G11").Value = '168'
    $ws.Range("H11").Value = $false
    $ws.Range("I11").Value = 'https://show.bilibili.com/platform/detail.html?id=80551'
    $ws.Range("J11").Value = '//i0.hdslb.com/bfs/openplatform/202401/MSS7qIQp1704695420767.jpeg'

    # Row 12
    $ws.Range("B12").Value = '2024-02-05'
    $ws.Range("C12").Value = '合肥·第十二届次元之门动漫游戏博览会-吴磊专场'
    $ws.Range("D12").Value = '经开区繁华大道与莲花路交叉口 百乐门大剧院'
    $ws.Range("E12").Value = '2024.02.05 09:00-02.05 17:00'
    $ws.Range("F12").Value = 11
    $ws.Range("G12").Value = '70'
    $ws.Range("H12").Value = $false
    $ws.Range("I12").Value = 'https://show.bilibili.com/platform/detail.html?id=81146'
    $ws.Range("J12").Value = '//i2.hdslb.com/bfs/openplatform/202401/QkgtYncY1705656564257.jpeg'

    # Row 13
    $ws.Range("B13").Value = '2024-02-13'
    $ws.Range("C13").Value = '合肥·环形宇宙动漫游戏嘉年华'
    $ws.Range("D13").Value = '山西路与太原路交叉口 挥动体育'
    $ws.Range("E13").Value = '2024.02.13 09:30-02.14 16:00'
    $ws.Range("F13").Value = 1722
    $ws.Range("G13").Value = '39'
    $ws.Range("H13").Value = $false
    $ws.Range("I13").Value = 'https://show.bilibili.com/platform/detail.html?id=80584'
    $ws.Range("J13").Value = '//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg'

    # Row 14
    $ws.Range("B14").Value = '2024-02-14'
    $ws.Range("C14").Value = '合肥·安徽马娘only'
    $ws.Range("D14").Value = '阜阳路16号 银瑞林国际大酒店'
    $ws.Range("E14").Value = '2024.02.14 10:00-02.14 17:00'
    $ws.Range("F14").Value = 58
    $ws.Range("G14").Value = '60'
    $ws.Range("H14").Value = $false
    $ws.Range("I14").Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
    $ws.Range("J14").Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'

    # Row 15
    $ws.Range("B15").Value = '2024-02-17'
    $ws.Range("C15").Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
    $ws.Range("D15").Value = '锦绣大道与清潭路交口东北角 李宁体育公园'
    $ws.Range("E15").Value = '2024.02.17 09:00-02.17 17:00'
    $ws.Range("F15").Value = 1101
    $ws.Range("G15").Value = '65'
    $ws.Range("H15").Value = $false
    $ws.Range("I15").Value = 'https://show.bilibili.com/platform/detail.html?id=79918'
    $ws.Range("J15").Value = '//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg'

    # Row 16
    $ws.Range("B16").Value = '2024-02-19'
    $ws.Range("C16").Value = '肥西·原神&崩铁&崩坏only'
    $ws.Range("D16").Value = '桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦'
    $ws.Range("E16").Value = '2024.02.19 09:00-02.19 17:00'
    $ws.Range("F16").Value = 271
    $ws.Range("G16").Value = '68'
    $ws.Range("H16").Value = $false
    $ws.Range("I16").Value = 'https://show.bilibili.com/platform/detail.html?id=78286'
    $ws.Range("J16").Value = '//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg'

    # Row 17
    $ws.Range("B17").Value = '2024-04-04'
    $ws.Range("C17").Value = '合肥·环形宇宙动漫游戏嘉年华—吴晛专场'
    $ws.Range("D17").Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
    $ws.Range("E17").Value = '2024.04.04 09:00-04.05 17:00'
    $ws.Range("F17").Value = 5486
    $ws.Range("G17").Value = '60'
    $ws.Range("H17").Value = $false
    $ws.Range("I17").Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
    $ws.Range("J17").Value = '//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg'

    # Remove the old last row (18), shrinking the table to 16 data rows
    $ws.Rows.Item(18).Delete()
}

Write-Output "edit complete"